# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.593.59'
$ws.Range('E2').Value = '  +5.22%  '
$ws.Range('D3').Value = '2.220.54'
$ws.Range('E3').Value = '  +3.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  -2.71%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.02'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.15%  '
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '2.549.84'
$ws.Range('E13').Value = '  +3.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.67'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.72'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.797'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.54'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.224.27'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('D19').Value = '41.526.75'
$ws.Range('E19').Value = '  +5.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.75'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('E21').Value = '  +5.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.05'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.99'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +10.40%  '
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.46'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.96'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.93'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.84%  '
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0622'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.33%  '
$ws.Range('E37').Value = '  -5.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.67'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000245'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +30.77%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.84'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E43').Value = '  +4.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.61'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0976'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.52'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.90%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').Value = '1.465.37'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.48'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.54%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.18'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +6.99%  '
